# Update the "Pais" sheet with refreshed COVID-19 country statistics
# and the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" banner text in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 12:05"

# Row 13 - India
$ws.Range("B13").Value = 139911
$ws.Range("C13").Value = 1375
$ws.Range("D13").Value = 57976
$ws.Range("E13").Value = 77896
$ws.Range("G13").Value = 15
$ws.Range("H13").Value = 4039

# Row 40 - Rumania
$ws.Range("B40").Value = 18283
$ws.Range("C40").Value = 213
$ws.Range("D40").Value = 11630
$ws.Range("E40").Value = 5460
$ws.Range("G40").Value = 8
$ws.Range("H40").Value = 1193

# Row 60 - Marruecos
$ws.Range("B60").Value = 7495
$ws.Range("C60").Value = 62
$ws.Range("D60").Value = 4737
$ws.Range("E60").Value = 2558
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 200

# Row 103 - Sri Lanka
$ws.Range("B103").Value = 1148
$ws.Range("C103").Value = 7
$ws.Range("E103").Value = 444

# Row 111 - Albania
$ws.Range("B111").Value = 1004
$ws.Range("C111").Value = 6
$ws.Range("D111").Value = 795
